$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The original data had two duplicate "ENVIRONMENT" rows (with stray
# leading/trailing spaces) at rows 3 and 5. Remove them, shifting the
# remaining rows (PROJECT, INSTRUMENT, ORGANISM part) up so the sheet
# ends up with a clean, deduplicated category ranking in A1:C5.

$ws.Rows("5").Delete()
$ws.Rows("3").Delete()
